# "fix some dir path"
# The FilePath column (F) for the DataNode rows referenced a
# "../../DataConfig/..." location that no longer exists; point it at the
# new "../resource/..." location instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F9").Value  = "../resource/Ini/Scene/1.xml"
$ws.Range("F10").Value = "../resource/Ini/Scene/2.xml"
$ws.Range("F11").Value = "../resource/Ini/Scene/3.xml"
$ws.Range("F12").Value = "../resource/Ini/Scene/4.xml"
$ws.Range("F13").Value = "../resource/Ini/Scene/5.xml"
$ws.Range("F14").Value = "../resource/Ini/Scene/6.xml"

# The sheet's active cell/selection moved as part of the same edit.
$ws.Range("F16").Select() | Out-Null
